$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matching original inlineStr values)
$textCells = @("D5","D6","D8","D10","D12","D13","D14","D19","D20","D21","D25","D26","D27","D28","D30","D31","D32","D35","D37","D38","D39","D40","D41","D43","D44","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "66.009.58"
$ws.Range("E2").Value = "  +6.25%  "
$ws.Range("D3").Value = "2.591.71"
$ws.Range("E3").Value = "  +6.20%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "587.89"
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "156.89"
$ws.Range("E6").Value = "  +8.29%  "
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("D8").Value = "0.545"
$ws.Range("E8").Value = "  +3.49%  "
$ws.Range("D9").Value = "2.620.82"
$ws.Range("E9").Value = "  +7.37%  "
$ws.Range("D10").Value = "0.116"
$ws.Range("E10").Value = "  +6.78%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("D13").Value = "0.362"
$ws.Range("E13").Value = "  +4.65%  "
$ws.Range("D14").Value = "29.62"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("E15").Value = "  +7.08%  "
$ws.Range("D16").Value = "3.065.12"
$ws.Range("E16").Value = "  +5.99%  "
$ws.Range("D17").Value = "65.998.66"
$ws.Range("E17").Value = "  +6.18%  "
$ws.Range("D18").Value = "2.625.58"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("D19").Value = "8.17"
$ws.Range("E19").Value = "  +6.31%  "
$ws.Range("D20").Value = "11.21"
$ws.Range("E20").Value = "  +5.08%  "
$ws.Range("D21").Value = "354.34"
$ws.Range("E21").Value = "  +10.88%  "
$ws.Range("E22").Value = "  +5.24%  "
$ws.Range("E23").Value = "  +6.34%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").Value = "10.10"
$ws.Range("E25").Value = "  +3.21%  "
$ws.Range("D26").Value = "66.24"
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").Value = "644.34"
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "0.0000107"
$ws.Range("E28").Value = "  +13.17%  "
$ws.Range("D29").Value = "2.719.56"
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("D30").Value = "1.52"
$ws.Range("E30").Value = "  +9.42%  "
$ws.Range("D31").Value = "0.993"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").Value = "8.27"
$ws.Range("E32").Value = "  +6.49%  "
$ws.Range("E33").Value = "  +6.15%  "
$ws.Range("E34").Value = "  +6.63%  "
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +11.45%  "
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "5.00"
$ws.Range("E37").Value = "  +8.63%  "
$ws.Range("D38").Value = "5.68"
$ws.Range("E38").Value = "  +9.31%  "
$ws.Range("D39").Value = "2.95"
$ws.Range("E39").Value = "  +10.45%  "
$ws.Range("D40").Value = "19.42"
$ws.Range("E40").Value = "  +5.90%  "
$ws.Range("D41").Value = "156.05"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("E42").Value = "  +3.66%  "
$ws.Range("D43").Value = "1.84"
$ws.Range("E43").Value = "  +9.30%  "
$ws.Range("D44").Value = "42.36"
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").Value = "0.0₆0318"
$ws.Range("E45").Value = "  +4.65%  "
$ws.Range("D46").Value = "163.29"
$ws.Range("E46").Value = "  +7.87%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "16.21"
$ws.Range("E48").Value = "  +6.06%  "
$ws.Range("D49").Value = "3.76"
$ws.Range("E49").Value = "  +7.35%  "
$ws.Range("D50").Value = "21.95"
$ws.Range("E50").Value = "  +10.68%  "
$ws.Range("D51").Value = "0.640"
$ws.Range("E51").Value = "  +6.54%  "
